$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

function Add-ItalicParagraphAfter($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.InsertParagraphAfter() | Out-Null
    $newp = $d.Paragraphs($paraIndex + 1)
    $rng = $newp.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $text
    $rng.Font.Italic = $true
}

# Insert from the bottom-most location upward so earlier paragraph indices stay valid.

# 3. After "Programa" body (paragraph 12) - English translation
Add-ItalicParagraphAfter 12 "Classification of manufacturing processes. Foundry. Powder metallurgy. Machining: processes, fundamentals and economic conditions. Machine tools. Mechanical conformation. Material joining processes. Computer Aided Manufacturing (CAM). Programming languages for numerical control. Numerical control machine tools. Product manufacturing sequence. Notions of automation of manufacturing processes. Rapid prototyping. Rapid prototyping systems (solid, liquid and powder)."

# 2. After "Programa resumido" body (paragraph 10) - English translation
Add-ItalicParagraphAfter 10 "Introduction to manufacturing processes. Material joining processes. Computer-aided design (CAD) review. Computer Aided Manufacturing (CAM). Flexible production systems. Rapid prototyping."

# 1. After "Objetivos" body (paragraph 6) - English translation
Add-ItalicParagraphAfter 6 "Develop theoretical and practical knowledge of the manufacturing processes of equipment and devices required for the development of products and prototypes. Know the requirements and effects of manufacturing processes in order to allow, interact, create and execute projects throughout your professional life."
